# Apply the Week's Games (FlashScore 2025-04-30) odds update
# Commit: "Atualizando o arquivo XLSX"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4
$ws.Range("G4").Value = 3.3
$ws.Range("I4").Value = 2.25
$ws.Range("T4").Value = 7.5
$ws.Range("W4").Value = 34
$ws.Range("AA4").Value = 6

# Row 5
$ws.Range("G5").Value = 2.02
$ws.Range("I5").Value = 3.55
$ws.Range("L5").Value = 1.36
$ws.Range("M5").Value = 2.67
$ws.Range("N5").Value = 2.05
$ws.Range("O5").Value = 1.6
$ws.Range("P5").Value = 1.47
$ws.Range("Q5").Value = 2.35
$ws.Range("R5").Value = 1.85
$ws.Range("S5").Value = 1.75
$ws.Range("T5").Value = 6.3
$ws.Range("U5").Value = 9
$ws.Range("V5").Value = 8.75
$ws.Range("W5").Value = 18
$ws.Range("X5").Value = 18
$ws.Range("Y5").Value = 32
$ws.Range("Z5").Value = 8
$ws.Range("AA5").Value = 6.2
$ws.Range("AB5").Value = 15.5
$ws.Range("AC5").Value = 80
$ws.Range("AD5").Value = 700
$ws.Range("AF5").Value = 18.5
$ws.Range("AH5").Value = 50
$ws.Range("AI5").Value = 35

# Row 7
$ws.Range("AD7").Value = 101

# Row 10
$ws.Range("H10").Value = 3.6
$ws.Range("I10").Value = 1.85
$ws.Range("J10").Value = 1.03
$ws.Range("L10").Value = 1.22
$ws.Range("N10").Value = 1.85
$ws.Range("O10").Value = 2
$ws.Range("R10").Value = 1.72
$ws.Range("T10").Value = 11
$ws.Range("U10").Value = 19
$ws.Range("V10").Value = 12
$ws.Range("AD10").Value = 201
$ws.Range("AE10").Value = 8
$ws.Range("AF10").Value = 9.5
$ws.Range("AH10").Value = 17

# Row 11
$ws.Range("L11").Value = 1.17
$ws.Range("M11").Value = 5
$ws.Range("N11").Value = 1.6
$ws.Range("O11").Value = 2.3

# Row 12
$ws.Range("G12").Value = 2.62
$ws.Range("H12").Value = 3.3
$ws.Range("I12").Value = 2.35
$ws.Range("L12").Value = 1.25
$ws.Range("M12").Value = 3.55
$ws.Range("S12").Value = 2.1
$ws.Range("U12").Value = 11.75
$ws.Range("V12").Value = 8.25
$ws.Range("W12").Value = 24
$ws.Range("X12").Value = 17
$ws.Range("Y12").Value = 21
$ws.Range("AA12").Value = 5.8
$ws.Range("AB12").Value = 10.5
$ws.Range("AD12").Value = 250
$ws.Range("AE12").Value = 7.7
$ws.Range("AF12").Value = 10.25
$ws.Range("AG12").Value = 7.8
$ws.Range("AH12").Value = 19.5
$ws.Range("AI12").Value = 15
$ws.Range("AJ12").Value = 20

# Row 15
$ws.Range("G15").Value = 2.32
$ws.Range("H15").Value = 2.87
$ws.Range("J15").Value = 1.14
$ws.Range("K15").Value = 5.3
$ws.Range("L15").Value = 1.57
$ws.Range("M15").Value = 2.3
$ws.Range("N15").Value = 2.7
$ws.Range("O15").Value = 1.42
$ws.Range("P15").Value = 1.62
$ws.Range("Q15").Value = 2.22
$ws.Range("R15").Value = 2.18
$ws.Range("S15").Value = 1.62
$ws.Range("T15").Value = 5.7
$ws.Range("V15").Value = 10.5
$ws.Range("X15").Value = 27
$ws.Range("Y15").Value = 55
$ws.Range("Z15").Value = 5.3
$ws.Range("AA15").Value = 6.1
$ws.Range("AB15").Value = 21
$ws.Range("AC15").Value = 150
$ws.Range("AE15").Value = 7.2
$ws.Range("AF15").Value = 17
$ws.Range("AG15").Value = 14
$ws.Range("AJ15").Value = 75

# Row 20
$ws.Range("J20").ClearContents()
$ws.Range("K20").ClearContents()
$ws.Range("L20").Value = 1.02

# Row 21
$ws.Range("G21").Value = 3.15
$ws.Range("H21").Value = 3.4
$ws.Range("K21").Value = 7.7
$ws.Range("M21").Value = 3.45
$ws.Range("P21").Value = 1.38
$ws.Range("Q21").Value = 2.8
$ws.Range("R21").Value = 1.65
$ws.Range("S21").Value = 2.12
$ws.Range("T21").Value = 10.75
$ws.Range("U21").Value = 17.5
$ws.Range("W21").Value = 40
$ws.Range("X21").Value = 26
$ws.Range("Y21").Value = 30
$ws.Range("Z21").Value = 7.7
$ws.Range("AA21").Value = 6.6
$ws.Range("AE21").Value = 8.25
$ws.Range("AF21").Value = 11
$ws.Range("AG21").Value = 8.75
$ws.Range("AI21").Value = 16.5
$ws.Range("AJ21").Value = 24

# Row 22
$ws.Range("G22").Value = 1.08
$ws.Range("H22").Value = 8.25
$ws.Range("I22").Value = 28
$ws.Range("L22").Value = 1.06
$ws.Range("M22").Value = 7.2
$ws.Range("N22").Value = 1.22
$ws.Range("O22").Value = 3.8
$ws.Range("P22").Value = 1.15
$ws.Range("Q22").Value = 4.65
$ws.Range("R22").Value = 1.98
$ws.Range("S22").Value = 1.75
$ws.Range("T22").Value = 14
$ws.Range("U22").Value = 8.5
$ws.Range("V22").Value = 12.5
$ws.Range("W22").Value = 7.3
$ws.Range("X22").Value = 10.75
$ws.Range("Y22").Value = 30
$ws.Range("Z22").Value = 30
$ws.Range("AA22").Value = 22
$ws.Range("AB22").Value = 37
$ws.Range("AC22").Value = 120
$ws.Range("AD22").Value = 700
$ws.Range("AE22").Value = 120
$ws.Range("AF22").Value = 500
$ws.Range("AG22").Value = 100
$ws.Range("AI22").Value = 600
$ws.Range("AJ22").Value = 250

# Row 23
$ws.Range("G23").Value = 3.6
$ws.Range("H23").Value = 3.65
$ws.Range("L23").Value = 1.18
$ws.Range("M23").Value = 4.3
$ws.Range("N23").Value = 1.55
$ws.Range("O23").Value = 2.3
$ws.Range("P23").Value = 1.3
$ws.Range("Q23").Value = 3.2
$ws.Range("S23").Value = 2.4
$ws.Range("T23").Value = 14.5
$ws.Range("U23").Value = 23
$ws.Range("X23").Value = 28
$ws.Range("Y23").Value = 28
$ws.Range("AA23").Value = 7.5
$ws.Range("AB23").Value = 11.75
$ws.Range("AC23").Value = 40
$ws.Range("AE23").Value = 10.25
$ws.Range("AF23").Value = 11.25

